{"js": "// [oldText, newText] pairs, in document order: the date line followed by the\n// 25 \"two-digit \u00f7 one-digit\" answer cells (diff regenerates all of them).\nconst replacements = [\n  [\"2025-06-06 Friday\", \"2025-06-07 Saturday\"],\n  [\"17\u00f76=2, 5\", \"50\u00f79=5, 5\"],\n  [\"92\u00f75=18, 2\", \"49\u00f76=8, 1\"],\n  [\"66\u00f73=22, 0\", \"23\u00f74=5, 3\"],\n  [\"95\u00f74=23, 3\", \"63\u00f72=31, 1\"],\n  [\"55\u00f75=11, 0\", \"48\u00f76=8, 0\"],\n  [\"91\u00f76=15, 1\", \"87\u00f73=29, 0\"],\n  [\"54\u00f73=18, 0\", \"95\u00f73=31, 2\"],\n  [\"49\u00f74=12, 1\", \"64\u00f72=32, 0\"],\n  [\"79\u00f77=11, 2\", \"93\u00f78=11, 5\"],\n  [\"84\u00f75=16, 4\", \"32\u00f79=3, 5\"],\n  [\"29\u00f76=4, 5\", \"81\u00f73=27, 0\"],\n  [\"37\u00f78=4, 5\", \"23\u00f79=2, 5\"],\n  [\"26\u00f75=5, 1\", \"64\u00f72=32, 0\"],\n  [\"11\u00f78=1, 3\", \"50\u00f77=7, 1\"],\n  [\"48\u00f73=16, 0\", \"74\u00f75=14, 4\"],\n  [\"55\u00f78=6, 7\", \"86\u00f78=10, 6\"],\n  [\"48\u00f77=6, 6\", \"65\u00f75=13, 0\"],\n  [\"36\u00f75=7, 1\", \"14\u00f72=7, 0\"],\n  [\"62\u00f74=15, 2\", \"13\u00f72=6, 1\"],\n  [\"15\u00f79=1, 6\", \"74\u00f73=24, 2\"],\n  [\"71\u00f77=10, 1\", \"85\u00f79=9, 4\"],\n  [\"90\u00f74=22, 2\", \"34\u00f78=4, 2\"],\n  [\"97\u00f75=19, 2\", \"91\u00f72=45, 1\"],\n  [\"80\u00f77=11, 3\", \"40\u00f79=4, 4\"],\n  [\"20\u00f75=4, 0\", \"75\u00f79=8, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# [oldText, newText] pairs, in document order: the date line followed by the\n# 25 \"two-digit \u00f7 one-digit\" answer cells (diff regenerates all of them).\n$replacements = @(\n  @(\"2025-06-06 Friday\", \"2025-06-07 Saturday\"),\n  @(\"17\u00f76=2, 5\", \"50\u00f79=5, 5\"),\n  @(\"92\u00f75=18, 2\", \"49\u00f76=8, 1\"),\n  @(\"66\u00f73=22, 0\", \"23\u00f74=5, 3\"),\n  @(\"95\u00f74=23, 3\", \"63\u00f72=31, 1\"),\n  @(\"55\u00f75=11, 0\", \"48\u00f76=8, 0\"),\n  @(\"91\u00f76=15, 1\", \"87\u00f73=29, 0\"),\n  @(\"54\u00f73=18, 0\", \"95\u00f73=31, 2\"),\n  @(\"49\u00f74=12, 1\", \"64\u00f72=32, 0\"),\n  @(\"79\u00f77=11, 2\", \"93\u00f78=11, 5\"),\n  @(\"84\u00f75=16, 4\", \"32\u00f79=3, 5\"),\n  @(\"29\u00f76=4, 5\", \"81\u00f73=27, 0\"),\n  @(\"37\u00f78=4, 5\", \"23\u00f79=2, 5\"),\n  @(\"26\u00f75=5, 1\", \"64\u00f72=32, 0\"),\n  @(\"11\u00f78=1, 3\", \"50\u00f77=7, 1\"),\n  @(\"48\u00f73=16, 0\", \"74\u00f75=14, 4\"),\n  @(\"55\u00f78=6, 7\", \"86\u00f78=10, 6\"),\n  @(\"48\u00f77=6, 6\", \"65\u00f75=13, 0\"),\n  @(\"36\u00f75=7, 1\", \"14\u00f72=7, 0\"),\n  @(\"62\u00f74=15, 2\", \"13\u00f72=6, 1\"),\n  @(\"15\u00f79=1, 6\", \"74\u00f73=24, 2\"),\n  @(\"71\u00f77=10, 1\", \"85\u00f79=9, 4\"),\n  @(\"90\u00f74=22, 2\", \"34\u00f78=4, 2\"),\n  @(\"97\u00f75=19, 2\", \"91\u00f72=45, 1\"),\n  @(\"80\u00f77=11, 3\", \"40\u00f79=4, 4\"),\n  @(\"20\u00f75=4, 0\", \"75\u00f79=8, 3\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute([ref]$oldText, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n"}
